$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add new row 79 (LASK vs Sturm Graz), copying formatting from row 78 first ---
$ws.Range("A78:V78").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "austria"
$ws.Range("C79").Value = "bundesliga"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45235.70833333334
$ws.Range("F79").Value = "LASK"
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = "Sturm Graz"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 2.73
$ws.Range("K79").Value = "29/10/2023 17:12"
$ws.Range("L79").Value = 2.87
$ws.Range("M79").Value = "05/11/2023 16:43"
$ws.Range("N79").Value = 3.56
$ws.Range("O79").Value = "29/10/2023 17:12"
$ws.Range("P79").Value = 3.3
$ws.Range("Q79").Value = "05/11/2023 16:43"
$ws.Range("R79").Value = 2.63
$ws.Range("S79").Value = "29/10/2023 17:12"
$ws.Range("T79").Value = 2.66
$ws.Range("U79").Value = "05/11/2023 16:43"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/austria/bundesliga/lask-linz-sturm-graz/GlUDEEbs/"

# --- Step 2: re-pair match rows whose betexplorer scrape order changed ---

# Row 44
$ws.Range("F44").Value = "Salzburg"
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = "BW Linz"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 1.16
$ws.Range("K44").Value = "16/09/2023 18:42"
$ws.Range("L44").Value = 1.26
$ws.Range("M44").Value = "23/09/2023 16:50"
$ws.Range("N44").Value = 8.36
$ws.Range("O44").Value = "16/09/2023 18:42"
$ws.Range("P44").Value = 6.59
$ws.Range("Q44").Value = "23/09/2023 16:57"
$ws.Range("R44").Value = 14.49
$ws.Range("S44").Value = "16/09/2023 18:42"
$ws.Range("T44").Value = 10.89
$ws.Range("U44").Value = "23/09/2023 16:57"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/austria/bundesliga/salzburg-bw-linz/6F2p9b54/"

# Row 45
$ws.Range("F45").Value = "A. Klagenfurt"
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = "A. Lustenau"
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1.68
$ws.Range("K45").Value = "17/09/2023 13:43"
$ws.Range("L45").Value = 1.69
$ws.Range("M45").Value = "23/09/2023 16:50"
$ws.Range("N45").Value = 4.12
$ws.Range("O45").Value = "17/09/2023 13:43"
$ws.Range("P45").Value = 4.04
$ws.Range("Q45").Value = "23/09/2023 16:50"
$ws.Range("R45").Value = 4.84
$ws.Range("S45").Value = "17/09/2023 13:43"
$ws.Range("T45").Value = 5.17
$ws.Range("U45").Value = "23/09/2023 16:50"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/austria/bundesliga/a-klagenfurt-a-lustenau/Eu6h7xzH/"

# Row 47
$ws.Range("F47").Value = "LASK"
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = "Hartberg"
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 1.59
$ws.Range("K47").Value = "17/09/2023 13:43"
$ws.Range("L47").Value = 1.74
$ws.Range("M47").Value = "24/09/2023 14:29"
$ws.Range("N47").Value = 4.48
$ws.Range("O47").Value = "17/09/2023 13:43"
$ws.Range("P47").Value = 4.17
$ws.Range("Q47").Value = "24/09/2023 14:29"
$ws.Range("R47").Value = 5.48
$ws.Range("S47").Value = "17/09/2023 13:43"
$ws.Range("T47").Value = 4.62
$ws.Range("U47").Value = "24/09/2023 14:14"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/austria/bundesliga/lask-linz-hartberg/0x2l8IKA/"

# Row 48
$ws.Range("F48").Value = "Altach"
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = "Austria Vienna"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 3.58
$ws.Range("K48").Value = "17/09/2023 13:43"
$ws.Range("L48").Value = 2.95
$ws.Range("M48").Value = "24/09/2023 14:29"
$ws.Range("N48").Value = 3.87
$ws.Range("O48").Value = "17/09/2023 13:43"
$ws.Range("P48").Value = 3.53
$ws.Range("Q48").Value = "24/09/2023 14:29"
$ws.Range("R48").Value = 1.99
$ws.Range("S48").Value = "17/09/2023 13:43"
$ws.Range("T48").Value = 2.47
$ws.Range("U48").Value = "24/09/2023 14:29"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/austria/bundesliga/altach-austria-vienna/zef15G4T/"

# Row 50
$ws.Range("F50").Value = "A. Lustenau"
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = "Salzburg"
$ws.Range("I50").Value = 4
$ws.Range("J50").Value = 9.88
$ws.Range("K50").Value = "23/09/2023 17:13"
$ws.Range("L50").Value = 13.38
$ws.Range("M50").Value = "30/09/2023 16:48"
$ws.Range("N50").Value = 6.64
$ws.Range("O50").Value = "23/09/2023 17:13"
$ws.Range("P50").Value = 7.85
$ws.Range("Q50").Value = "30/09/2023 16:48"
$ws.Range("R50").Value = 1.25
$ws.Range("S50").Value = "23/09/2023 17:13"
$ws.Range("T50").Value = 1.2
$ws.Range("U50").Value = "30/09/2023 16:34"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/austria/bundesliga/a-lustenau-salzburg/vstRLykA/"

# Row 51
$ws.Range("F51").Value = "BW Linz"
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = "A. Klagenfurt"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2.81
$ws.Range("K51").Value = "23/09/2023 17:13"
$ws.Range("L51").Value = 2.7
$ws.Range("M51").Value = "30/09/2023 16:50"
$ws.Range("N51").Value = 3.54
$ws.Range("O51").Value = "23/09/2023 17:13"
$ws.Range("P51").Value = 3.38
$ws.Range("Q51").Value = "30/09/2023 16:38"
$ws.Range("R51").Value = 2.5
$ws.Range("S51").Value = "23/09/2023 17:13"
$ws.Range("T51").Value = 2.77
$ws.Range("U51").Value = "30/09/2023 16:50"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/austria/bundesliga/bw-linz-a-klagenfurt/0jsNMHZ3/"

# Row 52
$ws.Range("F52").Value = "Wolfsberger AC"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "LASK"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 3.3
$ws.Range("K52").Value = "24/09/2023 13:42"
$ws.Range("L52").Value = 3.19
$ws.Range("M52").Value = "30/09/2023 16:43"
$ws.Range("N52").Value = 3.73
$ws.Range("O52").Value = "24/09/2023 13:42"
$ws.Range("P52").Value = 3.44
$ws.Range("Q52").Value = "30/09/2023 16:59"
$ws.Range("R52").Value = 2.13
$ws.Range("S52").Value = "24/09/2023 13:42"
$ws.Range("T52").Value = 2.35
$ws.Range("U52").Value = "30/09/2023 16:43"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/austria/bundesliga/wolfsberger-ac-lask-linz/6woXafSj/"

# Row 56
$ws.Range("F56").Value = "Austria Vienna"
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = "BW Linz"
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1.78
$ws.Range("K56").Value = "01/10/2023 16:12"
$ws.Range("L56").Value = 1.93
$ws.Range("M56").Value = "07/10/2023 16:58"
$ws.Range("N56").Value = 4.12
$ws.Range("O56").Value = "01/10/2023 16:12"
$ws.Range("P56").Value = 3.94
$ws.Range("Q56").Value = "07/10/2023 16:42"
$ws.Range("R56").Value = 4.22
$ws.Range("S56").Value = "01/10/2023 16:12"
$ws.Range("T56").Value = 3.81
$ws.Range("U56").Value = "07/10/2023 16:42"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/austria/bundesliga/austria-vienna-bw-linz/t817173e/"

# Row 57
$ws.Range("F57").Value = "A. Lustenau"
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = "Rapid Vienna"
$ws.Range("I57").Value = 5
$ws.Range("J57").Value = 5.42
$ws.Range("K57").Value = "01/10/2023 16:12"
$ws.Range("L57").Value = 7.27
$ws.Range("M57").Value = "07/10/2023 16:42"
$ws.Range("N57").Value = 4.45
$ws.Range("O57").Value = "01/10/2023 16:12"
$ws.Range("P57").Value = 5.21
$ws.Range("Q57").Value = "07/10/2023 16:42"
$ws.Range("R57").Value = 1.57
$ws.Range("S57").Value = "01/10/2023 16:12"
$ws.Range("T57").Value = 1.42
$ws.Range("U57").Value = "07/10/2023 16:42"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/austria/bundesliga/a-lustenau-rapid-vienna/2NkGaoY7/"

# Row 77
$ws.Range("F77").Value = "Hartberg"
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = "BW Linz"
$ws.Range("I77").Value = 2
$ws.Range("J77").Value = 1.92
$ws.Range("K77").Value = "29/10/2023 14:42"
$ws.Range("L77").Value = 1.81
$ws.Range("M77").Value = "05/11/2023 14:29"
$ws.Range("N77").Value = 3.69
$ws.Range("O77").Value = "29/10/2023 14:42"
$ws.Range("P77").Value = 3.71
$ws.Range("Q77").Value = "05/11/2023 14:29"
$ws.Range("R77").Value = 4.04
$ws.Range("S77").Value = "29/10/2023 14:42"
$ws.Range("T77").Value = 4.77
$ws.Range("U77").Value = "05/11/2023 14:29"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/austria/bundesliga/hartberg-bw-linz/pULUAWc6/"

# Row 78
$ws.Range("F78").Value = "Altach"
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = "Rapid Vienna"
$ws.Range("I78").Value = 2
$ws.Range("J78").Value = 3.49
$ws.Range("K78").Value = "29/10/2023 17:12"
$ws.Range("L78").Value = 4.46
$ws.Range("M78").Value = "05/11/2023 14:27"
$ws.Range("N78").Value = 3.8
$ws.Range("O78").Value = "29/10/2023 17:12"
$ws.Range("P78").Value = 3.78
$ws.Range("Q78").Value = "05/11/2023 14:27"
$ws.Range("R78").Value = 2.09
$ws.Range("S78").Value = "29/10/2023 17:12"
$ws.Range("T78").Value = 1.85
$ws.Range("U78").Value = "05/11/2023 14:27"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/austria/bundesliga/altach-rapid-vienna/6wNQBCr0/"
